$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new B24:B26 numbers the same look as the rest of column B
# (centered/wrapped, black font) by cloning the format already used by
# B23, rather than re-deriving each alignment/font property by hand.
$ws.Range("B23").Copy()
$ws.Range("B24:B26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new rows for the "Movimientos_de_Puntos_*" case, which was
# previously just a blank spacer row (B24) before the trailing blank row 27.
$ws.Range("A24").Value = "Movimientos_de_Puntos_MIX"
$ws.Range("B24").Value = 1162816939

$ws.Range("A25").Value = "Movimientos_de_Puntos_POS"
$ws.Range("B25").Value = 1145642605

$ws.Range("A26").Value = "Movimientos_de_Puntos_PRE"
$ws.Range("B26").Value = 1162676705

# Update view state to match the saved selection in the source file.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B24").Select()
